# Update row 2 (688571.SH) financial figures to the new reporting period.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: keep as text "001" (leading apostrophe forces text, matching
# the original inline-string/text storage instead of Excel auto-detecting a number).
$ws.Range("J2").Value = "'001"

# REPORT_DATE switches from 2020-06-30 to 2018-12-31 (stored as text, as before).
$ws.Range("N2").Value = "2018-12-31 00:00:00"

# Numeric financial figures / ratios.
$ws.Range("O2").Value = 1116067095.57
$ws.Range("P2").Value = 150714731.22
$ws.Range("Q2").Value = 291715724.81
$ws.Range("R2").Value = -12.427383843
$ws.Range("S2").Value = 255370172.84
$ws.Range("T2").Value = 3.1289325344
$ws.Range("U2").Value = 143778964.46
$ws.Range("V2").Value = -3.2544823115
$ws.Range("W2").Value = 297545549.2
$ws.Range("X2").Value = 226919998.06
$ws.Range("Y2").Value = -0.7313219379
$ws.Range("Z2").Value = 615728.1899999999
$ws.Range("AA2").Value = -61.1482598574
$ws.Range("AB2").Value = 818521546.37
$ws.Range("AC2").Value = 5.9043464767
$ws.Range("AD2").Value = 3.0399114802
$ws.Range("AE2").Value = -4.0958334824
$ws.Range("AF2").Value = 290.5918547554
$ws.Range("AG2").Value = 26.660184713
